# Fruta / hortaliza, semanal
# Insert 3 new weekly data rows at the top of the data block (row 17),
# pushing the existing rows 17-44 down to 20-47, then populate the 3
# new rows with this week's Espárragos prices for
# "Mapocho Venta Directa de Santiago".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 17 (Excel shifts
# everything below down and carries the row-above formatting, which is
# what keeps column D's date style (s="2") on the new rows).
$ws.Rows.Item(17).Resize(3).Insert()

# Row 17: Banquete
$ws.Cells.Item(17,1).Value2  = 12
$ws.Cells.Item(17,2).Value2  = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(17,3).Value2  = "Metropolitana"
$ws.Cells.Item(17,4).Value2  = 44483
$ws.Cells.Item(17,5).Value2  = 13
$ws.Cells.Item(17,6).Value2  = 300000000
$ws.Cells.Item(17,7).Value2  = "Espárragos"
$ws.Cells.Item(17,8).Value2  = "Sin especificar"
$ws.Cells.Item(17,9).Value2  = "Banquete"
$ws.Cells.Item(17,10).Value2 = 300
$ws.Cells.Item(17,11).Value2 = 1300
$ws.Cells.Item(17,12).Value2 = 1300
$ws.Cells.Item(17,13).Value2 = 1300
$ws.Cells.Item(17,14).Value2 = "$/kilo"
$ws.Cells.Item(17,15).Value2 = "Provincia de Linares"
$ws.Cells.Item(17,16).Value2 = 1300
$ws.Cells.Item(17,17).Value2 = 1
$ws.Cells.Item(17,18).Value2 = "Hortaliza"

# Row 18: Primera
$ws.Cells.Item(18,1).Value2  = 12
$ws.Cells.Item(18,2).Value2  = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(18,3).Value2  = "Metropolitana"
$ws.Cells.Item(18,4).Value2  = 44483
$ws.Cells.Item(18,5).Value2  = 13
$ws.Cells.Item(18,6).Value2  = 300000000
$ws.Cells.Item(18,7).Value2  = "Espárragos"
$ws.Cells.Item(18,8).Value2  = "Sin especificar"
$ws.Cells.Item(18,9).Value2  = "Primera"
$ws.Cells.Item(18,10).Value2 = 310
$ws.Cells.Item(18,11).Value2 = 1000
$ws.Cells.Item(18,12).Value2 = 1000
$ws.Cells.Item(18,13).Value2 = 1000
$ws.Cells.Item(18,14).Value2 = "$/kilo"
$ws.Cells.Item(18,15).Value2 = "Provincia de Linares"
$ws.Cells.Item(18,16).Value2 = 1000
$ws.Cells.Item(18,17).Value2 = 1
$ws.Cells.Item(18,18).Value2 = "Hortaliza"

# Row 19: Segunda
$ws.Cells.Item(19,1).Value2  = 12
$ws.Cells.Item(19,2).Value2  = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(19,3).Value2  = "Metropolitana"
$ws.Cells.Item(19,4).Value2  = 44483
$ws.Cells.Item(19,5).Value2  = 13
$ws.Cells.Item(19,6).Value2  = 300000000
$ws.Cells.Item(19,7).Value2  = "Espárragos"
$ws.Cells.Item(19,8).Value2  = "Sin especificar"
$ws.Cells.Item(19,9).Value2  = "Segunda"
$ws.Cells.Item(19,10).Value2 = 350
$ws.Cells.Item(19,11).Value2 = 700
$ws.Cells.Item(19,12).Value2 = 700
$ws.Cells.Item(19,13).Value2 = 700
$ws.Cells.Item(19,14).Value2 = "$/kilo"
$ws.Cells.Item(19,15).Value2 = "Provincia de Linares"
$ws.Cells.Item(19,16).Value2 = 700
$ws.Cells.Item(19,17).Value2 = 1
$ws.Cells.Item(19,18).Value2 = "Hortaliza"
